$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'63.797.52"
$ws.Cells.Item(2, 5).Value = '  +2.94%  '
$ws.Cells.Item(3, 4).Value = "'3.129.81"
$ws.Cells.Item(3, 5).Value = '  +1.69%  '
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).Value = "'589.34"
$ws.Cells.Item(5, 5).Value = '  +1.61%  '
$ws.Cells.Item(6, 4).Value = "'147.41"
$ws.Cells.Item(6, 5).Value = '  +3.65%  '
$ws.Cells.Item(7, 5).Value = '  -0.02%  '
$ws.Cells.Item(8, 4).Value = "'3.124.85"
$ws.Cells.Item(8, 5).Value = '  +1.79%  '
$ws.Cells.Item(9, 5).Value = '  +0.60%  '
$ws.Cells.Item(11, 5).Value = '  -0.22%  '
$ws.Cells.Item(12, 4).Value = "'0.470"
$ws.Cells.Item(12, 5).Value = '  +0.88%  '
$ws.Cells.Item(13, 4).Value = "'0.0000253"
$ws.Cells.Item(13, 5).Value = '  +4.93%  '
$ws.Cells.Item(14, 4).Value = "'37.24"
$ws.Cells.Item(14, 5).Value = '  +5.73%  '
$ws.Cells.Item(15, 5).Value = '  -0.67%  '
$ws.Cells.Item(16, 4).Value = "'3.645.66"
$ws.Cells.Item(16, 5).Value = '  +1.65%  '
$ws.Cells.Item(17, 5).Value = '  -1.43%  '
$ws.Cells.Item(18, 4).Value = "'63.658.50"
$ws.Cells.Item(19, 4).Value = "'3.127.48"
$ws.Cells.Item(19, 5).Value = '  +1.76%  '
$ws.Cells.Item(20, 4).Value = "'466.35"
$ws.Cells.Item(20, 5).Value = '  +4.06%  '
$ws.Cells.Item(21, 4).Value = "'14.36"
$ws.Cells.Item(21, 5).Value = '  +2.44%  '
$ws.Cells.Item(22, 5).Value = '  -0.24%  '
$ws.Cells.Item(23, 5).Value = '  +1.43%  '
$ws.Cells.Item(24, 4).Value = "'13.33"
$ws.Cells.Item(24, 5).Value = '  -3.03%  '
$ws.Cells.Item(25, 4).Value = "'82.30"
$ws.Cells.Item(25, 5).Value = '  +0.78%  '
$ws.Cells.Item(26, 5).Value = '  -0.05%  '
$ws.Cells.Item(27, 4).Value = "'9.00"
$ws.Cells.Item(27, 5).Value = '  +9.41%  '
$ws.Cells.Item(28, 5).Value = '  +1.85%  '
$ws.Cells.Item(29, 5).Value = '  -1.06%  '
$ws.Cells.Item(30, 5).Value = '  -0.04%  '
$ws.Cells.Item(31, 4).Value = "'6.89"
$ws.Cells.Item(31, 5).Value = '  +0.81%  '
$ws.Cells.Item(32, 4).Value = "'27.18"
$ws.Cells.Item(32, 5).Value = '  +1.19%  '
$ws.Cells.Item(33, 5).Value = '  -3.59%  '
$ws.Cells.Item(34, 4).Value = "'0.0₃0879"
$ws.Cells.Item(34, 5).Value = '  +10.48%  '
$ws.Cells.Item(35, 5).Value = '  +7.63%  '
$ws.Cells.Item(36, 5).Value = '  +16.72%  '
$ws.Cells.Item(37, 5).Value = '  +1.48%  '
$ws.Cells.Item(38, 4).Value = "'6.14"
$ws.Cells.Item(38, 5).Value = '  +1.43%  '
$ws.Cells.Item(39, 4).Value = "'459.10"
$ws.Cells.Item(39, 5).Value = '  +9.60%  '
$ws.Cells.Item(40, 4).Value = "'51.03"
$ws.Cells.Item(40, 5).Value = '  +1.70%  '
$ws.Cells.Item(41, 4).Value = "'8.75"
$ws.Cells.Item(41, 5).Value = '  -0.60%  '
$ws.Cells.Item(42, 4).Value = "'0.0375"
$ws.Cells.Item(42, 5).Value = '  +1.39%  '
$ws.Cells.Item(43, 4).Value = "'2.912.51"
$ws.Cells.Item(43, 5).Value = '  -0.38%  '
$ws.Cells.Item(44, 4).Value = "'0.279"
$ws.Cells.Item(44, 5).Value = '  +0.56%  '
$ws.Cells.Item(45, 5).Value = '  +2.34%  '
$ws.Cells.Item(46, 5).Value = '  +2.55%  '
$ws.Cells.Item(47, 4).Value = "'127.45"
$ws.Cells.Item(47, 5).Value = '  +3.00%  '
$ws.Cells.Item(48, 5).Value = '  +2.41%  '
$ws.Cells.Item(50, 5).Value = '  +0.39%  '
$ws.Cells.Item(51, 5).Value = '  +1.17%  '
